$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# New path values for the rac-geka images, keyed by the fitting/thread suffix
# (F = female, M = male, T = tetine/nozzle).
$pathF = "geka/rac-geka-F.png"
$pathM = "geka/rac-geka-M.png"
$pathT = "geka/rac-geka-T.png"

$ws.Range("A2").Value = $pathF
$ws.Range("A3").Value = $pathM
$ws.Range("A4").Value = $pathF
$ws.Range("A5").Value = $pathM
$ws.Range("A6").Value = $pathF
$ws.Range("A7").Value = $pathM
$ws.Range("A8").Value = $pathF
$ws.Range("A9").Value = $pathM
$ws.Range("A10").Value = $pathM
$ws.Range("A11").Value = $pathF
$ws.Range("A12").Value = $pathF
$ws.Range("A15").Value = $pathT
$ws.Range("A16").Value = $pathT
$ws.Range("A17").Value = $pathT
$ws.Range("A18").Value = $pathT
$ws.Range("A19").Value = $pathT

# Widen column A to fit the new, longer path text (bestFit/customWidth).
$ws.Columns.Item(1).ColumnWidth = 17

# Move the active selection to A19 like in the final saved file.
$ws.Range("A19").Select()
